# Script to calculate TGT resistor added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H9").Value = "9.27 k Ohms with 541 Ohm dampener"

$ws.Range("H10").Select()
